$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 17: F17 status -> Complete, fill in G17/H17 dates ---
$ws.Range("F17").Value = "Complete"
$ws.Range("G17").Value = (Get-Date -Year 2017 -Month 9 -Day 18)
$ws.Range("H17").Value = (Get-Date -Year 2017 -Month 9 -Day 18)
$ws.Range("G17").NumberFormat = "m/d/yyyy"
$ws.Range("H17").NumberFormat = "m/d/yyyy"

# --- Row 24: F24 status -> Complete, fill in G24/H24 dates ---
$ws.Range("F24").Value = "Complete"
$ws.Range("G24").Value = (Get-Date -Year 2017 -Month 9 -Day 18)
$ws.Range("H24").Value = (Get-Date -Year 2017 -Month 9 -Day 18)
$ws.Range("G24").NumberFormat = "m/d/yyyy"
$ws.Range("H24").NumberFormat = "m/d/yyyy"

# --- Row 28: align D28/E28 styles with the rest of the hidden data rows ---
$ws.Range("D28").HorizontalAlignment = -4108
$ws.Range("E28").WrapText = $true

# --- Update the view: scrolled position and active selection cell ---
$ws.Activate()
$appWin = $excel.ActiveWindow
$appWin.ScrollRow = 17
$ws.Range("E23").Select()
